$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (Ards and North Down Borough Council)
$ws.Range("I2").Value = 0
$ws.Range("N2").Value = 0.3035714285714285

# Row 3 (Derry City and Strabane District Council)
$ws.Range("G3").Value = 0.2857142857142857
$ws.Range("N3").Value = 0.2721428571428571

# Row 4 (Fermanagh and Omagh District Council)
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0.2701190476190476
